# pin分配.xlsx - update pinout doc
# - D4: "ENRF" -> "ENRF TIM3_CH1"
# - F4: "使能RF" -> "使能RF，可脉冲供电"
# - Fill in the sequence-number column (A2:A12) with 1..11
# - Widen column F
# - Move the active selection to D15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value() = "ENRF TIM3_CH1"
$ws.Range("F4").Value() = "使能RF，可脉冲供电"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value() = ($row - 1)
}

$ws.Columns.Item(6).ColumnWidth = 17.5

$ws.Range("D15").Select()
